$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) PODSUMOWANIE summary sheet: bump "last checked" timestamp for every
#    monitored profile from 21:47 to 21:51.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("PODSUMOWANIE")
for ($r = 2; $r -le 6; $r++) {
    $summary.Cells.Item($r, 2).Value = "2026-02-15 21:51"
}

# ---------------------------------------------------------------------------
# 2) Per-profile detail sheets: each got two fresh monitoring-run rows
#    appended (rows 6 and 7), logged at 2026-02-15 21:51. Row 6 re-uses the
#    "even" row style (same as row 2), row 7 re-uses the "odd" row style
#    (same as row 3 / row 5).
# ---------------------------------------------------------------------------
$details = @{
    "wszystkie-lublin" = @{ Count = 432; HasIds = $false; I6 = ""; I7 = ""; EvenRow = 2 }
    "artymiuk"         = @{ Count = 0;   HasIds = $false; I6 = ""; I7 = ""; EvenRow = 2 }
    "poqui"            = @{ Count = 5;   HasIds = $true;  I6 = "1951OR|17NeTz|17vbYq|18KAEc|183ger"; I7 = "183ger|18KAEc|17NeTz|17vbYq|1951OR"; EvenRow = 4 }
    "stylowepokoje"    = @{ Count = 2;   HasIds = $true;  I6 = "195dLc|16ZeYm"; I7 = "16ZeYm|195dLc"; EvenRow = 4 }
    "villahome"        = @{ Count = 0;   HasIds = $false; I6 = ""; I7 = ""; EvenRow = 2 }
}

foreach ($name in $details.Keys) {
    $info = $details[$name]
    $ws = $wb.Worksheets.Item($name)

    # Clone formatting (and, incidentally, the hidden column-I marker cell)
    # from existing plain "OK, nothing new" template rows onto the two new
    # rows before touching any values. Row 6 reuses whichever existing row
    # still carries the unhighlighted "even" row style (row 2, unless that
    # sheet's row 2 was the day the profile's first listings were found, in
    # which case row 4 is the plain one); row 7 always reuses row 5's
    # unhighlighted "odd" row style.
    $ws.Range("A$($info.EvenRow):I$($info.EvenRow)").Copy($ws.Range("A6:I6"))
    $ws.Range("A5:I5").Copy($ws.Range("A7:I7"))

    # Row 6
    $ws.Cells.Item(6, 1).Value = "2026-02-15 21:51"
    $ws.Cells.Item(6, 2).Value = $info.Count
    $ws.Cells.Item(6, 3).Value = 0
    $ws.Cells.Item(6, 4).Value = 0
    $ws.Cells.Item(6, 5).Value = 0
    $ws.Cells.Item(6, 6).Value = "—"
    $ws.Cells.Item(6, 7).Value = "—"
    $ws.Cells.Item(6, 8).Value = "OK"

    # Row 7
    $ws.Cells.Item(7, 1).Value = "2026-02-15 21:51"
    $ws.Cells.Item(7, 2).Value = $info.Count
    $ws.Cells.Item(7, 3).Value = 0
    $ws.Cells.Item(7, 4).Value = 0
    $ws.Cells.Item(7, 5).Value = 0
    $ws.Cells.Item(7, 6).Value = "—"
    $ws.Cells.Item(7, 7).Value = "—"
    $ws.Cells.Item(7, 8).Value = "OK"

    if ($info.HasIds) {
        # These profiles already track per-listing ids in column I; stamp
        # the freshly-shuffled id lists onto both new rows.
        $ws.Cells.Item(6, 9).Value = $info.I6
        $ws.Cells.Item(7, 9).Value = $info.I7
    } else {
        # These profiles never populate column I. Row 6's copy inherited an
        # empty I2, so drop it entirely; row 7's copy inherited an empty I3,
        # which is exactly the trailing empty marker we want to keep.
        $ws.Cells.Item(6, 9).Value = ""
    }

    # Row 5's trailing hidden-column marker moves to row 7 (the new last
    # row), so row 5 no longer carries it.
    $ws.Cells.Item(5, 9).Value = ""
}
